$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.715.52"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.513.77"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.63%  "
$ws.Range("D9").Value = "2.512.66"
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  +3.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("D14").Value = "2.980.27"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").Value = "69.627.89"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "2.517.23"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.50%  "
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D30").Value = "0.0₃0900"
$ws.Range("E30").Value = "  -2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "465.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  -12.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.01%  "
